$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fill in the remainder of existing row 16 (E16:L16) ---
$ws.Cells.Item(16, 5).Value = 17
$ws.Cells.Item(16, 6).Value = 7
$ws.Cells.Item(16, 7).Value = "Sunny, mild"
$ws.Cells.Item(16, 8).Value = $true
$ws.Cells.Item(16, 9).Value = "7 minutes 1 second"
$ws.Cells.Item(16, 10).Value = 421
$ws.Cells.Item(16, 11).Value = "Primary sweeps"
$ws.Cells.Item(16, 12).Value = "Worked downhill - but wind changed to downhill direction. Slightly trickier but Koda sourced after some redirecting."

# --- Row 17: copy cell formatting from row 16 (date / time styles) then set values ---
$ws.Range("A16:L16").Copy()
$ws.Range("A17:L17").PasteSpecial(-4122)

$ws.Cells.Item(17, 1).Value = 45808
$ws.Cells.Item(17, 2).Value = "PRESENCE"
$ws.Cells.Item(17, 3).Value = 0.54166666666666663
$ws.Cells.Item(17, 4).Value = 0.6875
$ws.Cells.Item(17, 5).Value = 15
$ws.Cells.Item(17, 6).Value = 7
$ws.Cells.Item(17, 7).Value = "Sunny, mild"
$ws.Cells.Item(17, 8).Value = $true
$ws.Cells.Item(17, 9).Value = "30 seconds"
$ws.Cells.Item(17, 10).Value = 30
$ws.Cells.Item(17, 11).Value = "Primary sweeps"
$ws.Cells.Item(17, 12).Value = "Worked uphill and found on the very first transect, it was so fast I nearly didn't belive it."

# --- Row 18: copy cell formatting from row 7 (ABSENCE / NA text / no-found pattern) ---
$ws.Range("A7:L7").Copy()
$ws.Range("A18:L18").PasteSpecial(-4122)

$ws.Cells.Item(18, 1).Value = 45809
$ws.Cells.Item(18, 2).Value = "ABSENCE"
$ws.Cells.Item(18, 3).Value = "NA"
$ws.Cells.Item(18, 4).Value = 0.45833333333333331
$ws.Cells.Item(18, 5).Value = 15
$ws.Cells.Item(18, 6).Value = 7
$ws.Cells.Item(18, 7).Value = "Sunny, mild"
$ws.Cells.Item(18, 8).Value = $false
$ws.Cells.Item(18, 9).Value = "13 minutes 23 seconds"
$ws.Cells.Item(18, 10).Value = 803
$ws.Cells.Item(18, 11).Value = "NA"
$ws.Cells.Item(18, 12).Value = "Worked uphill. Was pretty sure half way through there was no target because she was just loosely following me, but completed zig zags before finishing search."

# --- Row 19: copy cell formatting from row 17 (same pattern as row 16/17) ---
$ws.Range("A17:L17").Copy()
$ws.Range("A19:L19").PasteSpecial(-4122)

$ws.Cells.Item(19, 1).Value = 45809
$ws.Cells.Item(19, 2).Value = "PRESENCE"
$ws.Cells.Item(19, 3).Value = 0.59027777777777779
$ws.Cells.Item(19, 4).Value = 0.67708333333333337
$ws.Cells.Item(19, 5).Value = 14
$ws.Cells.Item(19, 6).Value = 9
$ws.Cells.Item(19, 7).Value = "Sunny, mild"
$ws.Cells.Item(19, 8).Value = $true
$ws.Cells.Item(19, 9).Value = "1 minute 20 seconds"
$ws.Cells.Item(19, 10).Value = 80
$ws.Cells.Item(19, 11).Value = "Primary sweeps"
$ws.Cells.Item(19, 12).Value = "Worked downhill. Found on the second sweep, nearly invisible under tussock."

# --- Update the view: scroll so column E is left-most visible, select I20 ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("I20").Select() | Out-Null
